# 17.6.1 — add 2021 (column M) and drop the stale "0.0" custom number
# format from the D:L data block (Excel collapses it back to General
# once the new column is filled in and the block is reformatted).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New column M: blank header-row cell (matches L2's border-only style) ----
$ws.Range("L2").Copy()
$ws.Range("M2").PasteSpecial(-4122)

# ---- New column M: year header 2021 (matches the other year header cells) ----
$ws.Range("F3").Copy()
$ws.Range("M3").PasteSpecial(-4122)
$ws.Range("M3").Value = 2021

# ---- Row 4 data: General number format across D4:M4, plus the new M4 value ----
$ws.Range("B4").Copy()
$ws.Range("D4:M4").PasteSpecial(-4122)
$ws.Range("M4").Value = 7105

# ---- Row 5 data: General number format across D5:M5, plus the new M5 value ----
$ws.Range("B5").Copy()
$ws.Range("D5:M5").PasteSpecial(-4122)
$ws.Range("M5").Value = 81079

# ---- Row 6 data (bottom-bordered row): General number format across D6:M6, plus new M6 value ----
$ws.Range("B6").Copy()
$ws.Range("D6:M6").PasteSpecial(-4122)
$ws.Range("M6").Value = 214139

$ws.Range("A1").Select()
